# Update "want to go" counts (column F) on the "展览" (Exhibition) sheet
# and the "全部类型" (All types) sheet, per the latest data refresh.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 67
$wsExhibit.Range("F4").Value = 168
$wsExhibit.Range("F6").Value = 5454
$wsExhibit.Range("F8").Value = 5405
$wsExhibit.Range("F10").Value = 9
$wsExhibit.Range("F11").Value = 1388
$wsExhibit.Range("F12").Value = 17
$wsExhibit.Range("F13").Value = 109

# --- Sheet "全部类型" ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 67
$wsAll.Range("F4").Value = 168
$wsAll.Range("F7").Value = 5454
$wsAll.Range("F9").Value = 5405
$wsAll.Range("F11").Value = 9
$wsAll.Range("F12").Value = 1388
$wsAll.Range("F13").Value = 17
$wsAll.Range("F14").Value = 109
